$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.904.80'
$ws.Range("E2").Value = '  +3.59%  '

$ws.Range("D3").Value = '2.534.32'
$ws.Range("E3").Value = '  +2.99%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''581.01'
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").Value = '''152.80'
$ws.Range("E6").Value = '  +4.01%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +1.18%  '

$ws.Range("D9").Value = '2.536.73'
$ws.Range("E9").Value = '  +3.04%  '

$ws.Range("E10").Value = '  +1.96%  '

$ws.Range("E11").Value = '  -1.77%  '

$ws.Range("D12").Value = '''5.29'
$ws.Range("E12").Value = '  +0.56%  '

$ws.Range("E13").Value = '  +0.69%  '

$ws.Range("E14").Value = '  +1.08%  '

$ws.Range("E15").Value = '  +2.86%  '

$ws.Range("D16").Value = '2.997.21'
$ws.Range("E16").Value = '  +3.00%  '

$ws.Range("D17").Value = '64.833.87'
$ws.Range("E17").Value = '  +3.75%  '

$ws.Range("D18").Value = '2.537.85'
$ws.Range("E18").Value = '  +3.07%  '

$ws.Range("D19").Value = '''8.03'
$ws.Range("E19").Value = '  +1.30%  '

$ws.Range("D20").Value = '''11.00'
$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("E21").Value = '  +3.57%  '

$ws.Range("D22").Value = '''329.94'
$ws.Range("E22").Value = '  +1.45%  '

$ws.Range("E23").Value = '  +2.95%  '

$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").Value = '''10.21'
$ws.Range("E25").Value = '  +2.03%  '

$ws.Range("D26").Value = '''65.87'

$ws.Range("D27").Value = '''643.83'
$ws.Range("E27").Value = '  +0.88%  '

$ws.Range("E28").Value = '  +7.95%  '

$ws.Range("E29").Value = '  +3.36%  '

$ws.Range("E30").Value = '  +5.20%  '

$ws.Range("D31").Value = '''0.995'
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("E32").Value = '  +2.04%  '

$ws.Range("E33").Value = '  +2.78%  '

$ws.Range("E34").Value = '  +3.68%  '

$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("E36").Value = '  +3.57%  '

$ws.Range("D37").Value = '''4.86'
$ws.Range("E37").Value = '  +2.65%  '

$ws.Range("D38").Value = '''5.65'
$ws.Range("E38").Value = '  +6.34%  '

$ws.Range("D39").Value = '''155.05'
$ws.Range("E39").Value = '  +3.07%  '

$ws.Range("E40").Value = '  +4.37%  '

$ws.Range("E41").Value = '  +1.40%  '

$ws.Range("D42").Value = '''18.93'
$ws.Range("E42").Value = '  +1.82%  '

$ws.Range("E43").Value = '  +5.97%  '

$ws.Range("D44").Value = '''161.63'
$ws.Range("E44").Value = '  +5.81%  '

$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("D46").Value = '0.0₆0301'
$ws.Range("E46").Value = '  -0.32%  '

$ws.Range("E47").Value = '  +2.34%  '

$ws.Range("D48").Value = '''3.66'
$ws.Range("E48").Value = '  +2.72%  '

$ws.Range("D49").Value = '''21.57'
$ws.Range("E49").Value = '  +6.43%  '

$ws.Range("E50").Value = '  +2.89%  '

$ws.Range("E51").Value = '  +2.27%  '
